$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 127 - another Christen Ford time entry on 12/7/2018
$ws.Range("A127").Value = "Christen Ford"
$ws.Range("B127").Value = 43441
$ws.Range("C127").Value = 0.375
$ws.Range("D127").Value = 0.45833333333333331
$ws.Range("E127").Formula = "=D127-C127"

# New row 128 - a second entry for the same day
$ws.Range("A128").Value = "Christen Ford"
$ws.Range("B128").Value = 43441
$ws.Range("C128").Value = 0.47916666666666669
$ws.Range("D128").Value = 0.59305555555555556
$ws.Range("E128").Formula = "=D128-C128"

# Copy the formatting (number formats / styles) of the last existing data
# row (126) down onto the two new rows so A127:E128 look like the rest of
# the log (text style for name, date style for B, time styles for C/D/E).
$ws.Range("A126:E126").Copy()
$ws.Range("A127:A128").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# The "Total Time Spent on Project" cell needs to sum through the new rows
$ws.Range("G3").Formula = "=SUM(E3:E128)"

# Make sure the summary formula recalculates against the extended range
$excel.CalculateFull()

# Reset the view: scroll back to the top and select G4 (matches the
# saved workbook state after finishing the last entries)
$null = $ws.Range("G4").Select()
